# "new init for scenario 29"
# The "Upper Left Cell" / "Lower Right Cell" index block for the Group
# Extraction section moved down one row (the extraction block that used to
# live at row 24 now lives at row 25), so the "Lower Right Cell" reference
# labels in D5:D11 need to shift from *24 to *25, and the active selection
# (previously the stale D5:D11 block) moves on to the next empty entry row,
# D12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Init")

$ws.Range("D5").Value  = "A25"
$ws.Range("D6").Value  = "B25"
$ws.Range("D7").Value  = "C25"
$ws.Range("D8").Value  = "G25"
$ws.Range("D9").Value  = "H25"
$ws.Range("D10").Value = "I25"
$ws.Range("D11").Value = "J25"

$ws.Range("D12").Select()
